# fmeasure_summary.xlsx - tempory save for 记忆
# Fill in the "CNN" Binkley results (C4:D13), compute their averages (C14:D14),
# tidy up the BT11 "CNN" column formatting (C15:C24), move the window position,
# and leave the selection where the author left it (K17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Binkley / CNN precision + recall values (previously empty) ---------
$ws.Range("C4").Value  = 0.78426839999999998
$ws.Range("D4").Value  = 0.78131220000000001
$ws.Range("C5").Value  = 0.78369962999999998
$ws.Range("D5").Value  = 0.77326010000000001
$ws.Range("C6").Value  = 0.78700124999999999
$ws.Range("D6").Value  = 0.80032175999999999
$ws.Range("C7").Value  = 0.75919705999999998
$ws.Range("D7").Value  = 0.77273583000000001
$ws.Range("C8").Value  = 0.78601109999999996
$ws.Range("D8").Value  = 0.79170214999999999
$ws.Range("C9").Value  = 0.7858619
$ws.Range("D9").Value  = 0.76981776999999996
$ws.Range("C10").Value = 0.79340560000000004
$ws.Range("D10").Value = 0.78021689999999999
$ws.Range("C11").Value = 0.77398299999999998
$ws.Range("D11").Value = 0.72584349999999997
$ws.Range("C12").Value = 0.81364460000000005
$ws.Range("D12").Value = 0.79706949999999999
$ws.Range("C13").Value = 0.78295890000000001
$ws.Range("D13").Value = 0.80299633999999998

# --- Averages for the newly filled columns -------------------------------
$ws.Range("C14").Formula = "=AVERAGE(C4:C13)"
$ws.Range("D14").Formula = "=AVERAGE(D4:D13)"

# --- Re-apply the bordered / centered look to the BT11 "CNN" column ------
# (distinct style entry from the untouched "without CNN" column next to it)
$bt11Cnn = $ws.Range("C15:C24")
$bt11Cnn.Borders.LineStyle = 1
$bt11Cnn.HorizontalAlignment = -4108
$bt11Cnn.VerticalAlignment = -4108
$bt11Cnn.Font.Name = "Times New Roman"
$bt11Cnn.Font.Size = 12
$bt11Cnn.Locked = $true

# --- Cosmetic window / selection state -----------------------------------
try { $excel.ActiveWindow.Left = 0 } catch {}
$ws.Range("K17").Select()
